# "Generate Report for Handback"
# Updates the localization-status report to reflect a completed handback:
#   - Overview sheet status text for both locale rows
#   - per-locale sheets (zh-cn, de-de) gain the "Latest Target File" hyperlink,
#     the "Latest Handback File" name and the "Latest Handback DateTime"
#   - column widths widen to fit the new, longer status text / handback links

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/324636875c39e7c047b78d5e9aa05b85ea20082a/e2e/"

$docOneName = "829a24b6-0fd2-4190-81f8-ca7ac1f38728.md"
$docTwoName = "baef03b3-719c-4d44-beb2-3883fd8dffc6.md"
$docOneUrl = $ghBase + $docOneName
$docTwoUrl = $ghBase + $docTwoName

# ---------------------------------------------------------------------------
# Overview sheet: both locale-status cells move from "Ready for handoff" to
# "Handed back: in sync with en-US" for both tracked documents.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

# Rebuild the hyperlinks in cell order (A2, I2, A3, I3) so the new
# "Latest Target File" links land between the existing "Source File Name"
# links, matching the handback report layout.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $docOneUrl, "", "", $docOneName)
$zh.Hyperlinks.Add($zh.Range("I2"), $docOneUrl, "", "", $docOneName)
$zh.Hyperlinks.Add($zh.Range("A3"), $docTwoUrl, "", "", $docTwoName)
$zh.Hyperlinks.Add($zh.Range("I3"), $docTwoUrl, "", "", $docTwoName)
$zh.Range("I2").Style = "HyperLink"
$zh.Range("I3").Style = "HyperLink"

$zh.Range("J2").Value = "829a24b6-0fd2-4190-81f8-ca7ac1f38728.09de99f5133a4972a2be7ed901b4040b2b979ee1.zh-cn.xlf"
$zh.Range("J3").Value = "baef03b3-719c-4d44-beb2-3883fd8dffc6.9bf79851f8aca950b05b7267021401f6867181aa.zh-cn.xlf"

$zh.Range("K2").Value = "2016-08-16 10:29:21"
$zh.Range("K3").Value = "2016-08-16 10:29:21"

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $docOneUrl, "", "", $docOneName)
$de.Hyperlinks.Add($de.Range("I2"), $docOneUrl, "", "", $docOneName)
$de.Hyperlinks.Add($de.Range("A3"), $docTwoUrl, "", "", $docTwoName)
$de.Hyperlinks.Add($de.Range("I3"), $docTwoUrl, "", "", $docTwoName)
$de.Range("I2").Style = "HyperLink"
$de.Range("I3").Style = "HyperLink"

$de.Range("J2").Value = "829a24b6-0fd2-4190-81f8-ca7ac1f38728.09de99f5133a4972a2be7ed901b4040b2b979ee1.de-de.xlf"
$de.Range("J3").Value = "baef03b3-719c-4d44-beb2-3883fd8dffc6.9bf79851f8aca950b05b7267021401f6867181aa.de-de.xlf"

$de.Range("K2").Value = "2016-08-16 10:29:29"
$de.Range("K3").Value = "2016-08-16 10:29:29"

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
